# Applies the scheduled-runner value updates captured in the commit diff.
# Each block below corresponds to one changed row in the "Unicorn Profits" sheets
# (columns H..N: currentAveragePrice*, LevePrice*, LeveProfit* computed figures).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 4647268.5
$ws.Range("I32").Value = 859
$ws.Range("J32").Value = 5808871
$ws.Range("K32").Value = 859
$ws.Range("L32").Value = 5808871
$ws.Range("M32").Value = -533
$ws.Range("N32").Value = -5809523

# Row 33
$ws.Range("H33").Value = 425.90625
$ws.Range("I33").Value = 405.35715
$ws.Range("J33").Value = 569.75
$ws.Range("K33").Value = 405.35715
$ws.Range("L33").Value = 569.75
$ws.Range("M33").Value = -176.35715
$ws.Range("N33").Value = -1027.75

# Row 100
$ws.Range("H100").Value = 4707.9287
$ws.Range("I100").Value = 2120.7144
$ws.Range("J100").Value = 5570.3335
$ws.Range("K100").Value = 2120.7144
$ws.Range("L100").Value = 5570.3335
$ws.Range("M100").Value = -1579.7144
$ws.Range("N100").Value = -6652.3335

# Row 121
$ws.Range("H121").Value = 4282.778
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 4718.125
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 14154.375
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -17648.375

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 22000
$ws.Range("J24").Value = 22000
$ws.Range("L24").Value = 22000
$ws.Range("N24").Value = -22748

# Row 61
$ws.Range("H61").Value = 276257.94
$ws.Range("I61").Value = 214326.31
$ws.Range("J61").Value = 388211.22
$ws.Range("K61").Value = 214326.31
$ws.Range("L61").Value = 388211.22
$ws.Range("M61").Value = -214114.31
$ws.Range("N61").Value = -388635.22

# Row 74
$ws.Range("H74").Value = 357375.34
$ws.Range("I74").Value = 770978.7
$ws.Range("K74").Value = 770978.7
$ws.Range("M74").Value = -770104.7

# Row 77
$ws.Range("H77").Value = 357375.34
$ws.Range("I77").Value = 770978.7
$ws.Range("K77").Value = 3854893.5
$ws.Range("M77").Value = -3850525.5

# Row 92
$ws.Range("H92").Value = 32309.4
$ws.Range("J92").Value = 32309.4
$ws.Range("L92").Value = 32309.4
$ws.Range("N92").Value = -37301.4

# Row 96
$ws.Range("H96").Value = 28215.5
$ws.Range("J96").Value = 28215.5
$ws.Range("L96").Value = 28215.5
$ws.Range("N96").Value = -33707.5

# Row 100
$ws.Range("H100").Value = 22000
$ws.Range("J100").Value = 22000
$ws.Range("L100").Value = 22000
$ws.Range("N100").Value = -24164

# Row 136
$ws.Range("H136").Value = 276257.94
$ws.Range("I136").Value = 214326.31
$ws.Range("J136").Value = 388211.22
$ws.Range("K136").Value = 642978.9299999999
$ws.Range("L136").Value = 1164633.66
$ws.Range("M136").Value = -640428.9299999999
$ws.Range("N136").Value = -1169733.66

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 3066.6667
$ws.Range("I107").Value = 3066.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3066.6667
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1146.6667
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4197.225
$ws.Range("I58").Value = 7537.467
$ws.Range("K58").Value = 7537.467
$ws.Range("M58").Value = -7334.467

# Row 94
$ws.Range("H94").Value = 4389.5
$ws.Range("I94").Value = 979.53845
$ws.Range("J94").Value = 8419.454
$ws.Range("K94").Value = 979.53845
$ws.Range("L94").Value = 8419.454
$ws.Range("M94").Value = -528.53845
$ws.Range("N94").Value = -9321.454

# Row 99
$ws.Range("H99").Value = 92822.37
$ws.Range("I99").Value = 144357.42
$ws.Range("J99").Value = 2636
$ws.Range("K99").Value = 144357.42
$ws.Range("L99").Value = 2636
$ws.Range("M99").Value = -142859.42
$ws.Range("N99").Value = -5632

# Row 105
$ws.Range("H105").Value = 857.2414
$ws.Range("I105").Value = 853.5714
$ws.Range("J105").Value = 960
$ws.Range("K105").Value = 853.5714
$ws.Range("L105").Value = 960
$ws.Range("M105").Value = 893.4286
$ws.Range("N105").Value = -4454

# Row 126
$ws.Range("H126").Value = 92822.37
$ws.Range("I126").Value = 144357.42
$ws.Range("J126").Value = 2636
$ws.Range("K126").Value = 433072.26
$ws.Range("L126").Value = 7908
$ws.Range("M126").Value = -430602.26
$ws.Range("N126").Value = -12848

# Row 136
$ws.Range("H136").Value = 4197.225
$ws.Range("I136").Value = 7537.467
$ws.Range("K136").Value = 22612.401
$ws.Range("M136").Value = -20062.401

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 583.8333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 583.8333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1751.4999
$ws.Range("N80").Value = -3623.4999
$ws.Range("M80").ClearContents()

# Row 83
$ws.Range("H83").Value = 583.8333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 583.8333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 5254.4997
$ws.Range("N83").Value = -14614.4997
$ws.Range("M83").ClearContents()

# Row 86
$ws.Range("H86").Value = 640.05884
$ws.Range("I86").Value = 676
$ws.Range("K86").Value = 2028
$ws.Range("M86").Value = -842

# Row 89
$ws.Range("H89").Value = 640.05884
$ws.Range("I89").Value = 676
$ws.Range("K89").Value = 6084
$ws.Range("M89").Value = -156

# Row 97
$ws.Range("H97").Value = 502
$ws.Range("J97").Value = 603
$ws.Range("L97").Value = 1809
$ws.Range("N97").Value = -2801

# Row 107
$ws.Range("H107").Value = 624.5833
$ws.Range("I107").Value = 247.66667
$ws.Range("J107").Value = 850.73334
$ws.Range("K107").Value = 743.00001
$ws.Range("L107").Value = 2552.20002
$ws.Range("M107").Value = 1176.99999
$ws.Range("N107").Value = -6392.20002

# Row 121
$ws.Range("H121").Value = 2276.889
$ws.Range("I121").Value = 825
$ws.Range("J121").Value = 3002.8333
$ws.Range("K121").Value = 2475
$ws.Range("L121").Value = 9008.499899999999
$ws.Range("M121").Value = -1165
$ws.Range("N121").Value = -11628.4999

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 1764.4783
$ws.Range("I126").Value = 1579.9445
$ws.Range("J126").Value = 2428.8
$ws.Range("K126").Value = 4739.833500000001
$ws.Range("L126").Value = 7286.400000000001
$ws.Range("M126").Value = -2269.833500000001
$ws.Range("N126").Value = -12226.4

# Row 136
$ws.Range("H136").Value = 20739.715
$ws.Range("J136").Value = 20739.715
$ws.Range("L136").Value = 62219.145
$ws.Range("N136").Value = -67319.145

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 550
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -5640
